$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C11").Value = 45174
